$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per the diff
$ws.Range("B2").Value = 558
$ws.Range("B3").Value = 444

# Remove row 4 entirely (A4, B4) so the sheet's used range becomes A1:B3
$ws.Range("A4:B4").Delete()
